$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 837.2560999999999
$ws.Range("I15").Value = 837.2560999999999
$ws.Range("K15").Value = 2511.7683
$ws.Range("M15").Value = -2342.7683
$ws.Range("H51").Value = 3979.7407
$ws.Range("I51").Value = 1541.8
$ws.Range("J51").Value = 4533.8184
$ws.Range("K51").Value = 1541.8
$ws.Range("L51").Value = 4533.8184
$ws.Range("M51").Value = -1057.8
$ws.Range("N51").Value = -5501.8184
$ws.Range("H135").Value = 17859336
$ws.Range("I135").Value = 2033.1666
$ws.Range("K135").Value = 18298.4994
$ws.Range("M135").Value = -15763.4994
$ws.Range("H137").Value = 3293.8635
$ws.Range("I137").Value = 3584.5
$ws.Range("J137").Value = 1986
$ws.Range("K137").Value = 10753.5
$ws.Range("L137").Value = 5958
$ws.Range("M137").Value = -8203.5
$ws.Range("N137").Value = -11058

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1549.9375
$ws.Range("I2").Value = 1280.4
$ws.Range("J2").Value = 1999.1666
$ws.Range("K2").Value = 1280.4
$ws.Range("L2").Value = 1999.1666
$ws.Range("M2").Value = -1167.4
$ws.Range("N2").Value = -2225.1666
$ws.Range("H45").Value = 1750.4
$ws.Range("I45").Value = 1667.3334
$ws.Range("J45").Value = 1875
$ws.Range("K45").Value = 1667.3334
$ws.Range("L45").Value = 1875
$ws.Range("M45").Value = -1290.3334
$ws.Range("N45").Value = -2629
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H61").Value = 22730642
$ws.Range("I61").Value = 31252642
$ws.Range("K61").Value = 31252642
$ws.Range("M61").Value = -31252430
$ws.Range("H116").Value = 1549.9375
$ws.Range("I116").Value = 1280.4
$ws.Range("J116").Value = 1999.1666
$ws.Range("K116").Value = 1280.4
$ws.Range("L116").Value = 1999.1666
$ws.Range("M116").Value = 1013.6
$ws.Range("N116").Value = -6587.1666
$ws.Range("H136").Value = 22730642
$ws.Range("I136").Value = 31252642
$ws.Range("K136").Value = 93757926
$ws.Range("M136").Value = -93755376

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1549.9375
$ws.Range("I3").Value = 1280.4
$ws.Range("J3").Value = 1999.1666
$ws.Range("K3").Value = 1280.4
$ws.Range("L3").Value = 1999.1666
$ws.Range("M3").Value = -1166.4
$ws.Range("N3").Value = -2227.1666
$ws.Range("H15").Value = 315
$ws.Range("I15").Value = 315
$ws.Range("K15").Value = 315
$ws.Range("M15").Value = -88
$ws.Range("H99").Value = 4679
$ws.Range("I99").Value = 2841.3333
$ws.Range("J99").Value = 6057.25
$ws.Range("K99").Value = 2841.3333
$ws.Range("L99").Value = 6057.25
$ws.Range("M99").Value = -1343.3333
$ws.Range("N99").Value = -9053.25
$ws.Range("H105").Value = 2085.1333
$ws.Range("I105").Value = 2151
$ws.Range("K105").Value = 2151
$ws.Range("M105").Value = -404
$ws.Range("H112").Value = 89999
$ws.Range("J112").Value = 89999
$ws.Range("L112").Value = 89999
$ws.Range("N112").Value = -92953
$ws.Range("H134").Value = 2553.721
$ws.Range("I134").Value = 2385.634
$ws.Range("K134").Value = 7156.902
$ws.Range("M134").Value = -4621.902

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H42").Value = 11693.75
$ws.Range("J42").Value = 11693.75
$ws.Range("L42").Value = 11693.75
$ws.Range("N42").Value = -12879.75
$ws.Range("H105").Value = 8675.058999999999
$ws.Range("I105").Value = 2689.6924
$ws.Range("K105").Value = 2689.6924
$ws.Range("M105").Value = -942.6923999999999
$ws.Range("H132").Value = 2310.5757
$ws.Range("I132").Value = 1804.5186
$ws.Range("K132").Value = 5413.5558
$ws.Range("M132").Value = -2883.5558
$ws.Range("H134").Value = 1437.76
$ws.Range("I134").Value = 1145.0952
$ws.Range("K134").Value = 3435.2856
$ws.Range("M134").Value = -900.2856000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 128.5
$ws.Range("I17").Value = 138.83333
$ws.Range("J17").Value = 118.166664
$ws.Range("K17").Value = 416.49999
$ws.Range("L17").Value = 354.499992
$ws.Range("M17").Value = -247.49999
$ws.Range("N17").Value = -692.499992
$ws.Range("H37").Value = 199996.75
$ws.Range("J37").Value = 199996.75
$ws.Range("L37").Value = 599990.25
$ws.Range("N37").Value = -600214.25
$ws.Range("H121").Value = 1096
$ws.Range("I121").Value = 958.3333
$ws.Range("J121").Value = 1171.091
$ws.Range("K121").Value = 2874.9999
$ws.Range("L121").Value = 3513.273
$ws.Range("M121").Value = -1564.9999
$ws.Range("N121").Value = -6133.272999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H111").Value = 39724.75
$ws.Range("J111").Value = 29966.334
$ws.Range("L111").Value = 29966.334
$ws.Range("N111").Value = -36100.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3869.8262
$ws.Range("I7").Value = 3414.4
$ws.Range("J7").Value = 4723.75
$ws.Range("K7").Value = 3414.4
$ws.Range("L7").Value = 4723.75
$ws.Range("M7").Value = -3302.4
$ws.Range("N7").Value = -4947.75
$ws.Range("H48").Value = 31268
$ws.Range("I48").Value = 25041
$ws.Range("J48").Value = 37495
$ws.Range("K48").Value = 25041
$ws.Range("L48").Value = -38817
$ws.Range("H68").Value = 3430.5417
$ws.Range("I68").Value = 3012.5
$ws.Range("J68").Value = 5520.75
$ws.Range("K68").Value = 3012.5
$ws.Range("L68").Value = 5520.75
$ws.Range("M68").Value = -2263.5
$ws.Range("N68").Value = -7018.75
$ws.Range("H71").Value = 3430.5417
$ws.Range("I71").Value = 3012.5
$ws.Range("J71").Value = 5520.75
$ws.Range("K71").Value = 15062.5
$ws.Range("L71").Value = 27603.75
$ws.Range("M71").Value = -11318.5
$ws.Range("N71").Value = -35091.75
$ws.Range("H82").Value = 4109.5625
$ws.Range("I82").Value = 2965.1428
$ws.Range("J82").Value = 4999.6665
$ws.Range("K82").Value = 2965.1428
$ws.Range("L82").Value = 4999.6665
$ws.Range("M82").Value = -2604.1428
$ws.Range("N82").Value = -5721.6665
$ws.Range("H85").Value = 4109.5625
$ws.Range("I85").Value = 2965.1428
$ws.Range("J85").Value = 4999.6665
$ws.Range("K85").Value = 2965.1428
$ws.Range("L85").Value = 4999.6665
$ws.Range("M85").Value = -1717.1428
$ws.Range("N85").Value = -7495.6665
$ws.Range("H87").Value = 65833.336
$ws.Range("J87").Value = 63000
$ws.Range("L87").Value = 63000
$ws.Range("N87").Value = -65246
$ws.Range("H90").Value = 65833.336
$ws.Range("J90").Value = 63000
$ws.Range("L90").Value = 189000
$ws.Range("N90").Value = -200232
$ws.Range("H93").Value = 2499
$ws.Range("I93").Value = 2148.6
$ws.Range("K93").Value = 2148.6
$ws.Range("M93").Value = -900.5999999999999
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H126").Value = 3869.8262
$ws.Range("I126").Value = 3414.4
$ws.Range("J126").Value = 4723.75
$ws.Range("K126").Value = 10243.2
$ws.Range("L126").Value = 14171.25
$ws.Range("M126").Value = -7773.200000000001
$ws.Range("N126").Value = -19111.25
$ws.Range("H132").Value = 100001816
$ws.Range("J132").Value = 250002000
$ws.Range("L132").Value = 750006000
$ws.Range("N132").Value = -750011060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 37999.5
$ws.Range("J121").Value = 37999.5
$ws.Range("L121").Value = 37999.5
$ws.Range("N121").Value = -41493.5
$ws.Range("H132").Value = 3833.1128
$ws.Range("I132").Value = 3970.8333
$ws.Range("K132").Value = 11912.4999
$ws.Range("M132").Value = -9382.499899999999
